# Update data obtained 2016-04-06: the "No aplicable" entry moves ahead of
# the "Equipamientos de salud..." entry in the tipo-local mapping table.
# The shared-string table gains "No aplicable" / its URL before the
# "Equipamientos de salud..." pair, which (since the rest of the sheet is
# untouched) surfaces as row 1 and row 2 swapping their displayed text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "No aplicable"
$ws.Range("B1").Value = "http://opendata.aragon.es/kos/iaest/tipo-local/no-aplicable"

$ws.Range("A2").Value = "Equipamientos de salud(ambulatorio, centro de salud, hospital,...)"
$ws.Range("B2").Value = "http://opendata.aragon.es/kos/iaest/tipo-local/equipamientos-de-saludambulatorio-centro-de-salud-hospital"
